# Update the Gdnf-Ret NATMI output (Ligand=Gdnf, Receptor=Ret) with the
# refreshed TPM-based numbers.
# - The "ECs" sending-cluster rows are dropped entirely.
# - The old "MuSCs" sending-cluster rows (rows 8-10) are removed, and the
#   table shrinks from 10 to 7 rows (header + 6 data rows): FAPs-as-sender
#   (rows 2-4) and MuSCs-as-sender (rows 5-7), each against Target cluster
#   ECs/FAPs/MuSCs, all with new expression/specificity values.
# Net effect: delete rows 8:10, then overwrite A2:T7 with the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old MuSCs-sender rows entirely (rows 8,9,10)
$ws.Range("A8:T10").EntireRow.Delete()

$data = @(
  @("FAPs","Gdnf","Ret","ECs",   1,0.3333333333333333,0.2660856666666667,0.798257,0.0824022153898117,0.1187122470109165,3,1,1.219350333333333,3.658051,0.2422674834150417,0.2576244469655636,0.3244516463452222,2.920064817107,0.0199633773503139,0.03058317698422676),
  @("FAPs","Gdnf","Ret","FAPs",  1,0.3333333333333333,0.2660856666666667,0.798257,0.0824022153898117,0.1187122470109165,3,1,2.913661333333334,8.740984000000001,0.5789028628226193,0.6155986258624717,0.7752835183208889,6.977551664888001,0.04770287839208809,0.07307909613296652),
  @("FAPs","Gdnf","Ret","MuSCs", 1,0.3333333333333333,0.2660856666666667,0.798257,0.0824022153898117,0.1187122470109165,2,1,0.9000629999999999,1.800126,0.178829653762339,0.1267769271719646,0.239493863397,1.436963180382,0.01473595964740971,0.01504997389372323),
  @("MuSCs","Gdnf","Ret","ECs",  2,1,2.9630225,5.926045,0.9175977846101883,0.8812877529890836,3,1,1.219350333333333,3.658051,0.2422674834150417,0.2576244469655636,3.612962473049167,21.677774838295,0.2223041060647278,0.2270412699813369),
  @("MuSCs","Gdnf","Ret","FAPs", 2,1,2.9630225,5.926045,0.9175977846101883,0.8812877529890836,3,1,2.913661333333334,8.740984000000001,0.5789028628226193,0.6155986258624717,8.633244088046668,51.79946452828001,0.5311999844305312,0.5425195297295052),
  @("MuSCs","Gdnf","Ret","MuSCs",2,1,2.9630225,5.926045,0.9175977846101883,0.8812877529890836,2,1,0.9000629999999999,1.800126,0.178829653762339,0.1267769271719646,2.6669069204175,10.66762768167,0.1640936941149293,0.1117269532782414)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
